# Implement basket-based elective scheduling with common time slots across all branches.
# Update the Section_A and Section_B timetable sheets with the new elective slot layout.

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Section_A")
$wsB = $wb.Worksheets.Item("Section_B")

# ---- Section_A : rows 2-3 (09:00-10:30, 10:30-12:00) ----
$dataA1 = New-Object 'object[,]' 2,5
$dataA1[0,0]="Free";  $dataA1[0,1]="CS304"; $dataA1[0,2]="HS101"; $dataA1[0,3]="CS303"; $dataA1[0,4]="CS309"
$dataA1[1,0]="Free";  $dataA1[1,1]="Free";  $dataA1[1,2]="CS309"; $dataA1[1,3]="Free";  $dataA1[1,4]="HS101"
$wsA.Range("B2:F3").Value = $dataA1

# ---- Section_A : rows 5-8 (13:00-14:30, 14:30-15:30, 15:30-17:00, 17:00-18:00) ----
$dataA2 = New-Object 'object[,]' 4,5
$dataA2[0,0]="Free";  $dataA2[0,1]="CS303"; $dataA2[0,2]="CS303"; $dataA2[0,3]="CS309"; $dataA2[0,4]="Free"
$dataA2[1,0]="Free";  $dataA2[1,1]="CS309 (Tutorial)"; $dataA2[1,2]="Free"; $dataA2[1,3]="Free"; $dataA2[1,4]="CS304 (Tutorial)"
$dataA2[2,0]="CS304"; $dataA2[2,1]="Free";  $dataA2[2,2]="Free";  $dataA2[2,3]="CS304"; $dataA2[2,4]="Free"
$dataA2[3,0]="Free";  $dataA2[3,1]="Free";  $dataA2[3,2]="Free";  $dataA2[3,3]="CS303 (Tutorial)"; $dataA2[3,4]="Free"
$wsA.Range("B5:F8").Value = $dataA2

# ---- Section_B : rows 2-3 (09:00-10:30, 10:30-12:00) ----
$dataB1 = New-Object 'object[,]' 2,5
$dataB1[0,0]="Free"; $dataB1[0,1]="CS309"; $dataB1[0,2]="HS101"; $dataB1[0,3]="Free"; $dataB1[0,4]="Free"
$dataB1[1,0]="Free"; $dataB1[1,1]="Free";  $dataB1[1,2]="Free";  $dataB1[1,3]="CS309"; $dataB1[1,4]="Free"
$wsB.Range("B2:F3").Value = $dataB1

# ---- Section_B : rows 5-8 (13:00-14:30, 14:30-15:30, 15:30-17:00, 17:00-18:00) ----
$dataB2 = New-Object 'object[,]' 4,5
$dataB2[0,0]="Free"; $dataB2[0,1]="CS303"; $dataB2[0,2]="CS304"; $dataB2[0,3]="CS304"; $dataB2[0,4]="CS304"
$dataB2[1,0]="Free"; $dataB2[1,1]="CS303 (Tutorial)"; $dataB2[1,2]="Free"; $dataB2[1,3]="CS304 (Tutorial)"; $dataB2[1,4]="Free"
$dataB2[2,0]="Free"; $dataB2[2,1]="HS101"; $dataB2[2,2]="CS303"; $dataB2[2,3]="CS303"; $dataB2[2,4]="CS309"
$dataB2[3,0]="Free"; $dataB2[3,1]="Free";  $dataB2[3,2]="Free";  $dataB2[3,3]="Free";  $dataB2[3,4]="CS309 (Tutorial)"
$wsB.Range("B5:F8").Value = $dataB2
